# Weekly refresh of the price series: insert a new latest-week row after
# row 29 (so it becomes row 30), shifting all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 30, pushing old rows 30-40 to 31-41.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new week's record.
$ws.Range("A30").Value = 4
$ws.Range("B30").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C30").Value = "Los Lagos"
$ws.Range("D30").Value = 44806
$ws.Range("E30").Value = 10
$ws.Range("F30").Value = 100112012
$ws.Range("G30").Value = "Espinaca"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 30
$ws.Range("K30").Value = 12000
$ws.Range("L30").Value = 12000
$ws.Range("M30").Value = 12000
$ws.Range("N30").Value = "$/cuna 10 kilos"
$ws.Range("O30").Value = "Región Metropolitana"
$ws.Range("P30").Value = 1200
$ws.Range("Q30").Value = 10
$ws.Range("R30").Value = "Hortaliza"
